# aggiornamento fino a 1/09/2021
# Extend the daily-data table (Sheet1) with 9 more rows (358-366),
# continuing the series from row 357 (date serial 44431 / 2021-08-23)
# through date serial 44440 (2021-09-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (incl. the date cell style used on column A) from the
# last existing data row down onto the new rows before filling in values,
# so the new rows pick up the same cell style (s="2" on column A) as the
# rest of the table.
$srcFormat = $ws.Range("A357:D357")
$dstFormat = $ws.Range("A358:D366")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)

# New data rows: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$data = @(
    @(358, 44432, 5, 11, 101.3357899585444),
    @(359, 44433, 0, 11, 101.3357899585444),
    @(360, 44434, 0, 11, 101.3357899585444),
    @(361, 44435, 5, 15, 138.1851681252879),
    @(362, 44436, 2, 17, 156.6098572086596),
    @(363, 44437, 1, 13, 119.7604790419162),
    @(364, 44438, 0, 13, 119.7604790419162),
    @(365, 44439, 4, 12, 110.5481345002303),
    @(366, 44440, 0, 12, 110.5481345002303)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
